# Applies the "Pooh Points: normal 20260221" update to the
# "Today_PoohPoints_SEC_ByOwner_2026-02-21" workbook:
#  - game clock moves from 7:49 to 0:26 of the 1st half, so every
#    player stat-line on the Players sheet is refreshed
#  - a box score correction reshuffles a few Undrafted bench names
#    (rows 14/16/17/18) and adds a new trailing row for Eduardo Klafke
#  - OwnerTotals starter totals are recalculated for the owners whose
#    starters stat-lines changed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Players")

# Row 2: James Scott
$ws.Cells.Item(2, 7).Value = "0:26 - 1st Half"

# Row 3: Thomas Haugh
$ws.Cells.Item(3, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(3, 8).Value = 17
$ws.Cells.Item(3, 9).Value = 14
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 16).Value = 17
$ws.Cells.Item(3, 17).Value = 4
$ws.Cells.Item(3, 18).Value = 7
$ws.Cells.Item(3, 19).Value = 2
$ws.Cells.Item(3, 20).Value = 4
$ws.Cells.Item(3, 21).Value = 4
$ws.Cells.Item(3, 22).Value = 4

# Row 4: Malik Dia
$ws.Cells.Item(4, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(4, 9).Value = 11
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 14).Value = 1
$ws.Cells.Item(4, 16).Value = 18
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = 11
$ws.Cells.Item(4, 21).Value = 1
$ws.Cells.Item(4, 22).Value = 3

# Row 5: Micah Handlogten
$ws.Cells.Item(5, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 6

# Row 6: Ilias Kamardine
$ws.Cells.Item(6, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(6, 8).Value = 0
$ws.Cells.Item(6, 9).Value = 6
$ws.Cells.Item(6, 14).Value = 2
$ws.Cells.Item(6, 15).Value = 1
$ws.Cells.Item(6, 16).Value = 17
$ws.Cells.Item(6, 18).Value = 5
$ws.Cells.Item(6, 21).Value = 4
$ws.Cells.Item(6, 22).Value = 4

# Row 7: Alex Condon
$ws.Cells.Item(7, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(7, 8).Value = 19
$ws.Cells.Item(7, 9).Value = 15
$ws.Cells.Item(7, 16).Value = 15
$ws.Cells.Item(7, 17).Value = 5
$ws.Cells.Item(7, 18).Value = 7

# Row 8: AJ Storr
$ws.Cells.Item(8, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(8, 8).Value = 2
$ws.Cells.Item(8, 9).Value = 4
$ws.Cells.Item(8, 10).Value = 2
$ws.Cells.Item(8, 16).Value = 11
$ws.Cells.Item(8, 17).Value = 2
$ws.Cells.Item(8, 18).Value = 6

# Row 9: Boogie Fland
$ws.Cells.Item(9, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(9, 8).Value = 3
$ws.Cells.Item(9, 9).Value = 2
$ws.Cells.Item(9, 10).Value = 3
$ws.Cells.Item(9, 12).Value = 2
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 16).Value = 13
$ws.Cells.Item(9, 18).Value = 3
$ws.Cells.Item(9, 20).Value = 2
$ws.Cells.Item(9, 21).Value = 2
$ws.Cells.Item(9, 22).Value = 2

# Row 10: Xaivian Lee
$ws.Cells.Item(10, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(10, 8).Value = 1
$ws.Cells.Item(10, 9).Value = 2
$ws.Cells.Item(10, 10).Value = 1
$ws.Cells.Item(10, 11).Value = 4
$ws.Cells.Item(10, 14).Value = 3
$ws.Cells.Item(10, 15).Value = 1
$ws.Cells.Item(10, 16).Value = 16
$ws.Cells.Item(10, 17).Value = 1
$ws.Cells.Item(10, 18).Value = 4
$ws.Cells.Item(10, 20).Value = 2

# Row 11: Rueben Chinyelu
$ws.Cells.Item(11, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(11, 8).Value = 4
$ws.Cells.Item(11, 9).Value = 1
$ws.Cells.Item(11, 10).Value = 2
$ws.Cells.Item(11, 16).Value = 14
$ws.Cells.Item(11, 21).Value = 1
$ws.Cells.Item(11, 22).Value = 2

# Row 12: Patton Pinkins
$ws.Cells.Item(12, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(12, 8).Value = 8
$ws.Cells.Item(12, 9).Value = 6
$ws.Cells.Item(12, 15).Value = 1
$ws.Cells.Item(12, 16).Value = 18
$ws.Cells.Item(12, 17).Value = 3
$ws.Cells.Item(12, 18).Value = 3

# Row 13: Urban Klavzar
$ws.Cells.Item(13, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(13, 8).Value = 8
$ws.Cells.Item(13, 9).Value = 6
$ws.Cells.Item(13, 15).Value = 2
$ws.Cells.Item(13, 16).Value = 11
$ws.Cells.Item(13, 17).Value = 2
$ws.Cells.Item(13, 18).Value = 2
$ws.Cells.Item(13, 19).Value = 2
$ws.Cells.Item(13, 20).Value = 2

# Row 14: Augusto Cassiá -> Travis Perry
$ws.Cells.Item(14, 4).Value = "Travis Perry"
$ws.Cells.Item(14, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(14, 8).Value = 5
$ws.Cells.Item(14, 9).Value = 3
$ws.Cells.Item(14, 10).Value = 1
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 14).Value = 1
$ws.Cells.Item(14, 16).Value = 10
$ws.Cells.Item(14, 18).Value = 2
$ws.Cells.Item(14, 19).Value = 1
$ws.Cells.Item(14, 20).Value = 2

# Row 15: Corey Chest
$ws.Cells.Item(15, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(15, 8).Value = 3
$ws.Cells.Item(15, 11).Value = 1
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 2
$ws.Cells.Item(15, 15).Value = 1
$ws.Cells.Item(15, 16).Value = 9

# Row 16: Eduardo Klafke -> Augusto Cassiá
$ws.Cells.Item(16, 4).Value = "Augusto Cassiá"
$ws.Cells.Item(16, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(16, 8).Value = 1
$ws.Cells.Item(16, 9).Value = 2
$ws.Cells.Item(16, 16).Value = 5
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = 2

# Row 17: Travis Perry -> Isaiah Brown
$ws.Cells.Item(17, 4).Value = "Isaiah Brown"
$ws.Cells.Item(17, 5).Value = "FLA"
$ws.Cells.Item(17, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(17, 8).Value = 1
$ws.Cells.Item(17, 9).Value = 3
$ws.Cells.Item(17, 10).Value = 1
$ws.Cells.Item(17, 14).Value = 2
$ws.Cells.Item(17, 15).Value = 1
$ws.Cells.Item(17, 16).Value = 8
$ws.Cells.Item(17, 17).Value = 1
$ws.Cells.Item(17, 18).Value = 1
$ws.Cells.Item(17, 21).Value = 1
$ws.Cells.Item(17, 22).Value = 2

# Row 18: Isaiah Brown -> Koren Johnson
$ws.Cells.Item(18, 4).Value = "Koren Johnson"
$ws.Cells.Item(18, 5).Value = "MISS"
$ws.Cells.Item(18, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 0
$ws.Cells.Item(18, 16).Value = 3
$ws.Cells.Item(18, 21).Value = 0
$ws.Cells.Item(18, 22).Value = 0

# Row 19 (new): Eduardo Klafke, Undrafted, MISS
# Force text format first so the date-like string "2026-02-21"
# is not auto-converted into an Excel date serial value, then
# clear the temporary formatting so the cell keeps the default style
$ws.Cells.Item(19, 1).NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "2026-02-21"
$ws.Cells.Item(19, 1).ClearFormats()
$ws.Cells.Item(19, 2).Value = "Undrafted"
$ws.Cells.Item(19, 3).Value = "No"
$ws.Cells.Item(19, 4).Value = "Eduardo Klafke"
$ws.Cells.Item(19, 5).Value = "MISS"
$ws.Cells.Item(19, 6).Value = "FLA@MISS"
$ws.Cells.Item(19, 7).Value = "0:26 - 1st Half"
$ws.Cells.Item(19, 8).Value = -1
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 0
$ws.Cells.Item(19, 16).Value = 1
$ws.Cells.Item(19, 17).Value = 0
$ws.Cells.Item(19, 18).Value = 1
$ws.Cells.Item(19, 19).Value = 0
$ws.Cells.Item(19, 20).Value = 0
$ws.Cells.Item(19, 21).Value = 0
$ws.Cells.Item(19, 22).Value = 0

# OwnerTotals: starter_pooh_total recalculated for owners with updated starters
$ws2 = $wb.Worksheets.Item("OwnerTotals")
$ws2.Cells.Item(2, 2).Value = 21   # Three Dawg Nite
$ws2.Cells.Item(3, 2).Value = 17   # The Backslashers
$ws2.Cells.Item(4, 2).Value = 4    # G-Flop
